# Edit coupons.xlsx:
#  1) Remove the leading "كوبون " ("Coupon ") prefix from the title (column A)
#     for rows 26-39.
#  2) Update the saved sheet view (zoom/scroll/selection) to match the
#     author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 26; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = [string]$cell.Value2
    if ($current.StartsWith("كوبون ")) {
        $cell.Value2 = $current.Substring(6)
    }
}

# Update view state: top-left visible cell, zoom level, and active selection.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.Zoom = 230
$ws.Range("B41").Select()
